$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("H23")
$c.Value = "URGENT(긴급) 0, TOP(상단 고정) 1, GENERAL(일반) 2"
$len = $c.Value.Length
Write-Host "len=$len"
$chars = $c.Characters(1, 11)
$chars.Font.Name = "Arial"
$chars2 = $c.Characters(12, 1)
$chars2.Font.Name = "돋움"
$chars2.Font.Size = 10
